# Refresh the cryptocurrency Price / Volume(1h) columns (cryptos.xlsx,
# GitHub Actions run on Sat Sep 23 09:14:31 UTC 2023).
#
# All D/E cells in this sheet hold plain text (e.g. "1.598.81", "26.664.61")
# rather than numbers, so every D-column write below first forces the cell's
# number format to Text ("@") -- otherwise Excel's COM layer would silently
# reinterpret a value like "211.60" or "90.90" as the number 211.6 / 90.9 and
# drop the trailing zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.669.97"
$ws.Range("E2").Value = "  -0.12%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.599.11"
$ws.Range("E3").Value = "  +0.08%  "

# Row 4
$ws.Range("E4").Value = "  +0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "211.60"
$ws.Range("E5").Value = "  +0.17%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.515"
$ws.Range("E6").Value = "  +0.90%  "

# Row 7
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("E9").Value = "  +0.63%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.58"
$ws.Range("E10").Value = "  -0.88%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0838"
$ws.Range("E11").Value = "  -0.01%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.822.89"
$ws.Range("E12").Value = "  +0.05%  "

# Row 13
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.619.03"
$ws.Range("E13").Value = "  +1.31%  "

# Row 14
$ws.Range("B14").Value = "Polkadot"
$ws.Range("C14").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.03"
$ws.Range("E14").Value = "  +0.09%  "

# Row 15
$ws.Range("E15").Value = "  +0.26%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.23"
$ws.Range("E16").Value = "  -0.04%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.656.16"
$ws.Range("E17").Value = "  -0.18%  "

# Row 18
$ws.Range("E18").Value = "  +1.45%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "209.62"
$ws.Range("E19").Value = "  -0.23%  "

# Row 20
$ws.Range("E20").Value = "  +0.18%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.04"
$ws.Range("E21").Value = "  +4.44%  "

# Row 22
$ws.Range("E22").Value = "  +0.71%  "

# Row 23
$ws.Range("E23").Value = "  +0.84%  "

# Row 24
$ws.Range("E24").Value = "  +0.90%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.25"
$ws.Range("E25").Value = "  -0.90%  "

# Row 26
$ws.Range("E26").Value = "  +0.13%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.12"
$ws.Range("E27").Value = "  -0.86%  "

# Row 28
$ws.Range("E28").Value = "  -0.57%  "

# Row 29
$ws.Range("E29").Value = "  -0.11%  "

# Row 30
$ws.Range("E30").Value = "  +2.20%  "

# Row 31
$ws.Range("E31").Value = "  +0.37%  "

# Row 32
$ws.Range("E32").Value = "  +0.48%  "

# Row 33
$ws.Range("E33").Value = "  +1.51%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.285.63"
$ws.Range("E34").Value = "  -0.89%  "

# Row 35
$ws.Range("E35").Value = "  -6.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.46"
$ws.Range("E36").Value = "  +0.55%  "

# Row 37
$ws.Range("E37").Value = "  +1.17%  "

# Row 38
$ws.Range("E38").Value = "  -0.71%  "

# Row 39
$ws.Range("E39").Value = "  -0.81%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.06"
$ws.Range("E40").Value = "  +21.35%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.53"
$ws.Range("E41").Value = "  +2.82%  "

# Row 42
$ws.Range("E42").Value = "  +0.18%  "

# Row 43
$ws.Range("E43").Value = "  -0.78%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "64.01"
$ws.Range("E44").Value = "  +0.49%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.735.96"
$ws.Range("E45").Value = "  +0.10%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "90.90"
$ws.Range("E46").Value = "  +0.86%  "

# Row 47
$ws.Range("E47").Value = "  -3.22%  "

# Row 48
$ws.Range("E48").Value = "  +3.19%  "

# Row 49
$ws.Range("E49").Value = "  +0.57%  "

# Row 50
$ws.Range("E50").Value = "  -0.07%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.41"
$ws.Range("E51").Value = "  -1.20%  "

